$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = '42.890.23'
$ws.Range("E2").Value = '  -0.90%  '
$ws.Range("D3").Value = '2.369.46'
$ws.Range("E3").Value = '  -1.45%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '318.86'
$ws.Range("E5").Value = '  -2.94%  '
$ws.Range("D6").Value = '109.03'
$ws.Range("E6").Value = '  +2.59%  '
$ws.Range("D7").Value = '0.636'
$ws.Range("E7").Value = '  -2.84%  '
$ws.Range("E8").Value = '  +0.09%  '
$ws.Range("D9").Value = '0.625'
$ws.Range("E9").Value = '  -4.57%  '
$ws.Range("D10").Value = '42.10'
$ws.Range("E10").Value = '  -0.31%  '
$ws.Range("D11").Value = '0.0931'
$ws.Range("E11").Value = '  -1.19%  '
$ws.Range("D12").Value = '8.59'
$ws.Range("E12").Value = '  -2.00%  '
$ws.Range("E13").Value = '  -3.91%  '
$ws.Range("E14").Value = '  -0.12%  '
$ws.Range("E15").Value = '  -5.96%  '
$ws.Range("D16").Value = '2.726.94'
$ws.Range("E16").Value = '  -1.37%  '
$ws.Range("D17").Value = '2.386.22'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").Value = '42.878.17'
$ws.Range("E18").Value = '  -0.75%  '
$ws.Range("D19").Value = '7.73'
$ws.Range("E19").Value = '  -0.69%  '
$ws.Range("D21").Value = '76.31'
$ws.Range("D22").Value = '3.71'
$ws.Range("E22").Value = '  -1.13%  '
$ws.Range("D23").Value = '257.82'
$ws.Range("D24").Value = '2.35'
$ws.Range("E24").Value = '  -3.79%  '
$ws.Range("D25").Value = '9.48'
$ws.Range("E25").Value = '  -2.59%  '
$ws.Range("E26").Value = '  +0.12%  '
$ws.Range("E27").Value = '  -2.85%  '
$ws.Range("D28").Value = '23.09'
$ws.Range("E28").Value = '  -0.84%  '
$ws.Range("E29").Value = '  +2.59%  '
$ws.Range("D30").Value = '37.24'
$ws.Range("E30").Value = '  -0.28%  '
$ws.Range("D31").Value = '172.40'
$ws.Range("E31").Value = '  -2.28%  '
$ws.Range("D32").Value = '0.0897'
$ws.Range("E32").Value = '  -4.95%  '
$ws.Range("D33").Value = '6.08'
$ws.Range("E33").Value = '  +1.60%  '
$ws.Range("D34").Value = '2.96'
$ws.Range("E34").Value = '  -6.50%  '
$ws.Range("E35").Value = '  +11.75%  '
$ws.Range("E36").Value = '  -3.00%  '
$ws.Range("E37").Value = '  -3.96%  '
$ws.Range("D38").Value = '0.0366'
$ws.Range("E38").Value = '  +0.16%  '
$ws.Range("E39").Value = '  -4.48%  '
$ws.Range("D40").Value = '2.69'
$ws.Range("E40").Value = '  -5.73%  '
$ws.Range("E41").Value = '  +2.82%  '
$ws.Range("D42").Value = '1.51'
$ws.Range("E42").Value = '  -5.16%  '
$ws.Range("D43").Value = '71.94'
$ws.Range("E43").Value = '  +2.21%  '
$ws.Range("E44").Value = '  +0.12%  '
$ws.Range("D45").Value = '12.37'
$ws.Range("E45").Value = '  -0.21%  '
$ws.Range("D46").Value = '113.47'
$ws.Range("E46").Value = '  -7.92%  '
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("E48").Value = '  -0.74%  '
$ws.Range("D49").Value = '86.34'
$ws.Range("E49").Value = '  -5.14%  '
$ws.Range("D50").Value = '76.89'
$ws.Range("E50").Value = '  +6.41%  '
$ws.Range("E51").Value = '  -0.87%  '
